$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("G3").Value = 0.0107
$ws.Range("H3").Value = 0.0107

# Row 5
$ws.Range("G5").Value = 0.61
$ws.Range("H5").Value = 1.22

# Row 6
$ws.Range("G6").Value = 0.03374
$ws.Range("H6").Value = 0.3374

# Row 7
$ws.Range("G7").Value = 0.0077
$ws.Range("H7").Value = 0.077

# Row 9: Manufacturer Part Number 1 changed from numeric barcode to text part number
$ws.Range("C9").Value = "GRM1885C1H102JA01D"
$ws.Range("G9").Value = 0.0144
$ws.Range("H9").Value = 0.0144

# Row 11
$ws.Range("G11").Value = 0.0737
$ws.Range("H11").Value = 0.5159

# Row 12
$ws.Range("G12").Value = 0.0413
$ws.Range("H12").Value = 0.0413

# Row 14
$ws.Range("G14").Value = 0.0149
$ws.Range("H14").Value = 0.0149

# Row 18
$ws.Range("G18").Value = 0.189
$ws.Range("H18").Value = 0.189

# Row 20: clear out unit price / subtotal (no longer applicable)
$ws.Range("G20").ClearContents()
$ws.Range("H20").ClearContents()

# Row 22
$ws.Range("G22").Value = 0.301

# Row 30
$ws.Range("G30").Value = 0.019
$ws.Range("H30").Value = 0.19

# Row 33
$ws.Range("G33").Value = 0.0319
$ws.Range("H33").Value = 0.0638

# Row 38
$ws.Range("G38").Value = 0.0044
$ws.Range("H38").Value = 0.0044
